$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.203.53"
$ws.Range("E2").Value = "  -0.62%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.859.68"
$ws.Range("E3").Value = "  -1.28%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7041"
$ws.Range("E5").Value = "  -1.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.29"
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3105"
$ws.Range("E8").Value = "  -0.68%  "
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07773"
$ws.Range("E9").Value = "  -3.26%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.17"
$ws.Range("E10").Value = "  -4.50%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08007"
$ws.Range("E11").Value = "  -3.91%  "
$ws.Range("B12").Value = "Litecoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "93.44"
$ws.Range("E12").Value = "  +0.54%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.169"
$ws.Range("E13").Value = "  -1.51%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.827.05"
$ws.Range("E14").Value = "  -3.46%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6949"
$ws.Range("E15").Value = "  -3.64%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.353"
$ws.Range("E16").Value = "  +0.21%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "29.214.36"
$ws.Range("E17").Value = "  -0.63%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008284"
$ws.Range("E18").Value = "  -1.93%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "251.39"
$ws.Range("E19").Value = "  +4.26%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.127.37"
$ws.Range("E20").Value = "  -1.32%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.07"
$ws.Range("E21").Value = "  -1.41%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.497"
$ws.Range("E23").Value = "  -4.64%  "
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1549"
$ws.Range("E25").Value = "  -2.52%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.954"
$ws.Range("E26").Value = "  -1.04%  "
$ws.Range("E27").Value = "  -2.88%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.76"
$ws.Range("E28").Value = "  +1.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.498"
$ws.Range("E29").Value = "  -0.63%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.275"
$ws.Range("E30").Value = "  -3.29%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.254"
$ws.Range("E31").Value = "  -2.09%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.208"
$ws.Range("E32").Value = "  +0.76%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05247"
$ws.Range("E33").Value = "  -2.27%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.884"
$ws.Range("E34").Value = "  -3.44%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7423"
$ws.Range("E35").Value = "  -1.08%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.154"
$ws.Range("E36").Value = "  -2.53%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.704"
$ws.Range("E37").Value = "  +0.06%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01860"
$ws.Range("E38").Value = "  -1.53%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.251.67"
$ws.Range("E39").Value = "  -3.12%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.734"
$ws.Range("E40").Value = "  -0.51%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.259"
$ws.Range("E41").Value = "  -5.19%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8971"
$ws.Range("E42").Value = "  -2.48%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "110.87"
$ws.Range("E43").Value = "  -0.88%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "71.21"
$ws.Range("E44").Value = "  -3.92%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.001"
$ws.Range("E45").Value = "  +0.00%  "
$ws.Range("E46").Value = "  -0.34%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.024.35"
$ws.Range("E47").Value = "  -0.80%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5200"
$ws.Range("E48").Value = "  -0.36%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.775"
$ws.Range("E49").Value = "  -1.83%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.389"
$ws.Range("E50").Value = "  -1.26%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4284"
$ws.Range("E51").Value = "  -2.75%  "
